# "this added 15-1024 report"
# The Sep-24 P&L workbook's "Fuel" expense line (B42) increases from 7000 to
# 8000. Every other changed cell (F41, F43, B51) is a formula that depends on
# B42 (directly or via the Total Expense chain), so Excel recalculates them
# automatically once the input cell is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fuel (row 42, column B) : 7000 -> 8000
$ws.Range("B42").Value = 8000

# Leave the cursor where the author ended up after making the edit.
$ws.Range("B43").Select()
